$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 497.875
$ws.Range("I103").Value = 498.16666
$ws.Range("J103").Value = 497
$ws.Range("K103").Value = 1494.49998
$ws.Range("L103").Value = 1491
$ws.Range("M103").Value = -908.4999800000001
$ws.Range("N103").Value = -2663

$ws.Range("H106").Value = 1683.3334
$ws.Range("I106").Value = 1025
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 1025
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -394
$ws.Range("N106").Value = -4262

$ws.Range("H109").Value = 28842
$ws.Range("I109").Value = 23000
$ws.Range("J109").Value = 34684
$ws.Range("K109").Value = 23000
$ws.Range("L109").Value = 34684
$ws.Range("M109").Value = -21613
$ws.Range("N109").Value = -37458

$ws.Range("H129").Value = 599
$ws.Range("J129").Value = 1100
$ws.Range("L129").Value = 3300
$ws.Range("N129").Value = -13300

$ws.Range("H132").Value = 2179.4595
$ws.Range("I132").Value = 1672.9
$ws.Range("J132").Value = 4350.4287
$ws.Range("K132").Value = 5018.700000000001
$ws.Range("L132").Value = 13051.2861
$ws.Range("M132").Value = -2488.700000000001
$ws.Range("N132").Value = -18111.2861

$ws.Range("H135").Value = 337.92856
$ws.Range("I135").Value = 335.84616
$ws.Range("J135").Value = 365
$ws.Range("K135").Value = 3022.61544
$ws.Range("L135").Value = 3285
$ws.Range("M135").Value = -487.61544
$ws.Range("N135").Value = -8355

$ws.Range("H137").Value = 4738.5557
$ws.Range("I137").Value = 1591
$ws.Range("J137").Value = 6517.609
$ws.Range("K137").Value = 4773
$ws.Range("L137").Value = 19552.827
$ws.Range("M137").Value = -2223
$ws.Range("N137").Value = -24652.827

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 983.6667
$ws.Range("I4").Value = 1120.4
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 1120.4
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -1004.4
$ws.Range("N4").Value = -532

$ws.Range("H61").Value = 1538.8518
$ws.Range("I61").Value = 1117.3158
$ws.Range("J61").Value = 2540
$ws.Range("K61").Value = 1117.3158
$ws.Range("L61").Value = 2540
$ws.Range("M61").Value = -905.3158000000001
$ws.Range("N61").Value = -2964

$ws.Range("H136").Value = 1538.8518
$ws.Range("I136").Value = 1117.3158
$ws.Range("J136").Value = 2540
$ws.Range("K136").Value = 3351.9474
$ws.Range("L136").Value = 7620
$ws.Range("M136").Value = -801.9474
$ws.Range("N136").Value = -12720

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 392.55554
$ws.Range("I22").Value = 181.7
$ws.Range("J22").Value = 656.125
$ws.Range("K22").Value = 181.7
$ws.Range("L22").Value = 656.125
$ws.Range("M22").Value = -8.699999999999989
$ws.Range("N22").Value = -1002.125

$ws.Range("H105").Value = 1794
$ws.Range("I105").Value = 1353.875
$ws.Range("J105").Value = 2234.125
$ws.Range("K105").Value = 1353.875
$ws.Range("L105").Value = 2234.125
$ws.Range("M105").Value = 393.125
$ws.Range("N105").Value = -5728.125

$ws.Range("H134").Value = 1489.8823
$ws.Range("I134").Value = 960.46155
$ws.Range("J134").Value = 3210.5
$ws.Range("K134").Value = 2881.38465
$ws.Range("L134").Value = 9631.5
$ws.Range("M134").Value = -346.38465
$ws.Range("N134").Value = -14701.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 11584.444
$ws.Range("I22").Value = 17025
$ws.Range("K22").Value = 17025
$ws.Range("M22").Value = -16675

$ws.Range("H31").Value = 27874.545
$ws.Range("I31").Value = 2539.2942
$ws.Range("J31").Value = 43826.37
$ws.Range("K31").Value = 2539.2942
$ws.Range("L31").Value = 43826.37
$ws.Range("M31").Value = -2244.2942
$ws.Range("N31").Value = -44416.37

$ws.Range("H34").Value = 27874.545
$ws.Range("I34").Value = 2539.2942
$ws.Range("J34").Value = 43826.37
$ws.Range("K34").Value = 2539.2942
$ws.Range("L34").Value = 43826.37
$ws.Range("M34").Value = -2337.2942
$ws.Range("N34").Value = -44230.37

$ws.Range("H58").Value = 1094.4166
$ws.Range("I58").Value = 973.7273
$ws.Range("J58").Value = 1359.9333
$ws.Range("K58").Value = 973.7273
$ws.Range("L58").Value = 1359.9333
$ws.Range("M58").Value = -770.7273
$ws.Range("N58").Value = -1765.9333

$ws.Range("H118").Value = 25500
$ws.Range("J118").Value = 25500
$ws.Range("L118").Value = 25500
$ws.Range("N118").Value = -28814

$ws.Range("H136").Value = 1094.4166
$ws.Range("I136").Value = 973.7273
$ws.Range("J136").Value = 1359.9333
$ws.Range("K136").Value = 2921.1819
$ws.Range("L136").Value = 4079.7999
$ws.Range("M136").Value = -371.1819
$ws.Range("N136").Value = -9179.7999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5124.9116
$ws.Range("I132").Value = 5626.7036
$ws.Range("J132").Value = 3189.4285
$ws.Range("K132").Value = 16880.1108
$ws.Range("L132").Value = 9568.2855
$ws.Range("M132").Value = -14350.1108
$ws.Range("N132").Value = -14628.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 55
$ws.Range("I20").Value = 55
$ws.Range("K20").Value = 55
$ws.Range("M20").Value = 171

$ws.Range("H82").Value = 1270.5714
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1270.5714
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1270.5714
$ws.Range("N82").Value = -1992.5714
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 1270.5714
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1270.5714
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1270.5714
$ws.Range("N85").Value = -3766.5714
$ws.Range("M85").ClearContents()

$ws.Range("H132").Value = 3810.907
$ws.Range("I132").Value = 3999.6
$ws.Range("J132").Value = 2985.375
$ws.Range("K132").Value = 11998.8
$ws.Range("L132").Value = 8956.125
$ws.Range("M132").Value = -9468.799999999999
$ws.Range("N132").Value = -14016.125

$ws.Range("H136").Value = 4154.757
$ws.Range("I136").Value = 1452
$ws.Range("K136").Value = 4356
$ws.Range("M136").Value = -1806

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0

$ws.Range("H136").Value = 2286.0833
$ws.Range("I136").Value = 2087.1072
$ws.Range("J136").Value = 2982.5
$ws.Range("K136").Value = 6261.321599999999
$ws.Range("L136").Value = 8947.5
$ws.Range("M136").Value = -3711.321599999999
$ws.Range("N136").Value = -14047.5
